$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.186.87"
Set-TextValue $ws.Range("E2") "  +2.96%  "
Set-TextValue $ws.Range("D3") "1.895.15"
Set-TextValue $ws.Range("E3") "  +0.00%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  -0.22%  "
Set-TextValue $ws.Range("D5") "325.02"
Set-TextValue $ws.Range("E5") "  +3.31%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  -0.23%  "
Set-TextValue $ws.Range("D7") "0.5161"
Set-TextValue $ws.Range("E7") "  +0.25%  "
Set-TextValue $ws.Range("D8") "0.3989"
Set-TextValue $ws.Range("E8") "  +1.64%  "
Set-TextValue $ws.Range("D9") "0.08418"
Set-TextValue $ws.Range("E9") "  -0.14%  "
Set-TextValue $ws.Range("D10") "42.65"
Set-TextValue $ws.Range("E10") "  +0.62%  "
Set-TextValue $ws.Range("D11") "1.115"
Set-TextValue $ws.Range("E11") "  -0.11%  "
Set-TextValue $ws.Range("D12") "23.37"
Set-TextValue $ws.Range("E12") "  +12.91%  "
Set-TextValue $ws.Range("D13") "6.421"
Set-TextValue $ws.Range("E13") "  +2.48%  "
Set-TextValue $ws.Range("D14") "1.892.32"
Set-TextValue $ws.Range("E14") "  -0.19%  "
Set-TextValue $ws.Range("D15") "7.325"
Set-TextValue $ws.Range("E15") "  +0.44%  "
Set-TextValue $ws.Range("D16") "1.002"
Set-TextValue $ws.Range("E16") "  -0.20%  "
Set-TextValue $ws.Range("D17") "94.02"
Set-TextValue $ws.Range("E17") "  +0.86%  "
Set-TextValue $ws.Range("D18") "0.00001111"
Set-TextValue $ws.Range("E18") "  +0.66%  "
Set-TextValue $ws.Range("D19") "0.06645"
Set-TextValue $ws.Range("E19") "  -1.29%  "
Set-TextValue $ws.Range("D20") "18.21"
Set-TextValue $ws.Range("E20") "  +2.00%  "
Set-TextValue $ws.Range("E21") "  -0.21%  "
Set-TextValue $ws.Range("D22") "5.945"
Set-TextValue $ws.Range("E22") "  -1.15%  "
Set-TextValue $ws.Range("D23") "30.188.02"
Set-TextValue $ws.Range("E23") "  +2.88%  "
Set-TextValue $ws.Range("D24") "11.26"
Set-TextValue $ws.Range("E24") "  +1.06%  "
Set-TextValue $ws.Range("D25") "2.226"
Set-TextValue $ws.Range("E25") "  +0.51%  "
Set-TextValue $ws.Range("D26") "2.108.04"
Set-TextValue $ws.Range("E26") "  -0.18%  "
Set-TextValue $ws.Range("E27") "  +3.80%  "
Set-TextValue $ws.Range("D28") "161.46"
Set-TextValue $ws.Range("E28") "  +1.36%  "
Set-TextValue $ws.Range("D29") "2.357"
Set-TextValue $ws.Range("E29") "  -3.09%  "
Set-TextValue $ws.Range("D30") "129.13"
Set-TextValue $ws.Range("E30") "  +1.38%  "
Set-TextValue $ws.Range("D31") "1.095"
Set-TextValue $ws.Range("E31") "  +3.22%  "
Set-TextValue $ws.Range("D32") "0.1056"
Set-TextValue $ws.Range("E32") "  +0.60%  "
Set-TextValue $ws.Range("D33") "6.082"
Set-TextValue $ws.Range("E33") "  -0.61%  "
Set-TextValue $ws.Range("D34") "3.756"
Set-TextValue $ws.Range("E34") "  +2.67%  "
Set-TextValue $ws.Range("D35") "0.02500"
Set-TextValue $ws.Range("E35") "  +0.79%  "
Set-TextValue $ws.Range("D36") "0.06547"
Set-TextValue $ws.Range("E36") "  -0.46%  "
Set-TextValue $ws.Range("D37") "5.267"
Set-TextValue $ws.Range("E37") "  +1.43%  "
Set-TextValue $ws.Range("D38") "0.2199"
Set-TextValue $ws.Range("E38") "  +0.00%  "
Set-TextValue $ws.Range("D39") "1.217"
Set-TextValue $ws.Range("E39") "  -0.99%  "
Set-TextValue $ws.Range("D40") "11.74"
Set-TextValue $ws.Range("E40") "  +4.09%  "
Set-TextValue $ws.Range("D41") "0.6494"
Set-TextValue $ws.Range("E41") "  -0.40%  "
Set-TextValue $ws.Range("D42") "8.720"
Set-TextValue $ws.Range("E42") "  -3.63%  "
Set-TextValue $ws.Range("D43") "1.230"
Set-TextValue $ws.Range("E43") "  -0.26%  "
Set-TextValue $ws.Range("D44") "0.6098"
Set-TextValue $ws.Range("E44") "  +0.43%  "
Set-TextValue $ws.Range("D45") "13.16"
Set-TextValue $ws.Range("E45") "  -0.14%  "
Set-TextValue $ws.Range("E46") "  +0.86%  "
Set-TextValue $ws.Range("D47") "2.052"
Set-TextValue $ws.Range("E47") "  +0.05%  "
Set-TextValue $ws.Range("D48") "1.233"
Set-TextValue $ws.Range("E48") "  +0.37%  "
Set-TextValue $ws.Range("D49") "124.46"
Set-TextValue $ws.Range("E49") "  +0.81%  "
Set-TextValue $ws.Range("D51") "79.02"
Set-TextValue $ws.Range("E51") "  +1.69%  "
